# Update the dSF column (column F) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = -7
    "F3"  = 1
    "F4"  = -2
    "F6"  = 3
    "F11" = 1
    "F15" = -6
    "F17" = -3
    "F20" = -2
    "F22" = 0
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
